$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E = 5 (剩余 / remaining), Column F = 6 (开始时间 / start date)
$updates = @(
    @{ Row = 2; E = 12 }
    @{ Row = 3; E = 12 }
    @{ Row = 4; E = 12 }
    @{ Row = 5; E = 8 }
    @{ Row = 6; E = 12 }
    @{ Row = 7; E = 8 }
    @{ Row = 8; E = 12 }
    @{ Row = 9; E = 8 }
    @{ Row = 10; E = 5 }
    @{ Row = 11; E = 12 }
    @{ Row = 12; E = 8 }
    @{ Row = 13; E = 12 }
    @{ Row = 14; E = 12 }
    @{ Row = 15; E = 12 }
    @{ Row = 16; E = 2 }
    @{ Row = 17; E = 8 }
    @{ Row = 18; E = 1 }
    @{ Row = 19; E = 1 }
    @{ Row = 20; E = 1 }
    @{ Row = 21; E = 1 }
    @{ Row = 22; E = 8 }
    @{ Row = 23; E = 8 }
    @{ Row = 24; E = 8 }
    @{ Row = 25; E = 8 }
    @{ Row = 26; E = 8 }
    @{ Row = 27; E = 6 }
    @{ Row = 28; E = 1 }
    @{ Row = 29; E = 1 }
    @{ Row = 30; E = 1 }
    @{ Row = 31; E = 1 }
    @{ Row = 32; E = 1 }
    @{ Row = 33; E = 1 }
    @{ Row = 34; E = 1 }
    @{ Row = 35; E = 1 }
    @{ Row = 37; E = 1 }
    @{ Row = 38; E = 1 }
    @{ Row = 39; E = 1 }
    @{ Row = 40; E = 5 }
    @{ Row = 41; E = 5 }
    @{ Row = 42; E = 1 }
    @{ Row = 43; E = 8 }
    @{ Row = 44; E = 5 }
    @{ Row = 45; E = 8 }
    @{ Row = 46; E = 5 }
    @{ Row = 47; E = 1 }
    @{ Row = 48; E = 5 }
    @{ Row = 49; E = 6 }
    @{ Row = 50; E = 6 }
    @{ Row = 51; E = 6 }
    @{ Row = 52; E = 6 }
    @{ Row = 53; E = 6 }
    @{ Row = 54; E = 6 }
    @{ Row = 55; E = 6 }
    @{ Row = 56; E = 6 }
    @{ Row = 57; E = 6 }
    @{ Row = 58; E = 10; F = 20251217 }
    @{ Row = 59; E = 10; F = 20251217 }
    @{ Row = 60; E = 10; F = 20251217 }
    @{ Row = 61; E = 6 }
    @{ Row = 62; E = 10; F = 20251217 }
    @{ Row = 63; E = 10; F = 20251217 }
    @{ Row = 64; E = 10; F = 20251217 }
    @{ Row = 65; E = 1 }
    @{ Row = 66; E = 1 }
    @{ Row = 67; E = 1 }
    @{ Row = 68; E = 1 }
    @{ Row = 69; E = 1 }
    @{ Row = 70; E = 2 }
    @{ Row = 71; E = 2 }
    @{ Row = 72; E = 2 }
    @{ Row = 73; E = 2 }
    @{ Row = 74; E = 2 }
    @{ Row = 75; E = 2 }
    @{ Row = 76; E = 2 }
    @{ Row = 77; E = 5 }
    @{ Row = 78; E = 5 }
    @{ Row = 79; E = 5 }
    @{ Row = 80; E = 5 }
    @{ Row = 81; E = 5 }
    @{ Row = 82; E = 5 }
    @{ Row = 83; E = 5 }
    @{ Row = 84; E = 5 }
    @{ Row = 85; E = 5 }
    @{ Row = 86; E = 5 }
    @{ Row = 87; E = 5 }
    @{ Row = 88; E = 5 }
    @{ Row = 89; E = 5 }
    @{ Row = 90; E = 5 }
    @{ Row = 91; E = 8 }
    @{ Row = 92; E = 5 }
    @{ Row = 93; E = 5 }
    @{ Row = 94; E = 1 }
    @{ Row = 95; E = 4 }
    @{ Row = 96; E = 2 }
    @{ Row = 97; E = 2 }
    @{ Row = 98; E = 2 }
    @{ Row = 99; E = 2 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    if ($u.ContainsKey("F")) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
    }
}
